$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("scannedStatus") - mark scanned trees as Yes/No
$ws.Range("E2").Value = "Yes"
$ws.Range("E3").Value = "No"
$ws.Range("E4").Value = "Yes"
$ws.Range("E5").Value = "No"
$ws.Range("E6").Value = "Yes"
$ws.Range("E7").Value = "No"
$ws.Range("E8").Value = "No"
$ws.Range("E9").Value = "Yes"
$ws.Range("E10").Value = "No"
$ws.Range("E11").Value = "No"
$ws.Range("E12").Value = "No"
$ws.Range("E13").Value = "Yes"

# Column F ("Date") - rescan dates moved forward two days (45922 -> 45924)
$ws.Range("F3").Value = 45924
$ws.Range("F5").Value = 45924
$ws.Range("F7").Value = 45924
$ws.Range("F8").Value = 45924
$ws.Range("F10").Value = 45924
$ws.Range("F11").Value = 45924
$ws.Range("F12").Value = 45924
$ws.Range("F14").Value = 45924
$ws.Range("F15").Value = 45924
$ws.Range("F16").Value = 45924
$ws.Range("F17").Value = 45924
$ws.Range("F18").Value = 45924

# Rows 19-25 haven't been scanned yet - clear their dates
$ws.Range("F19").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("F21").ClearContents()
$ws.Range("F22").ClearContents()
$ws.Range("F23").ClearContents()
$ws.Range("F24").ClearContents()
$ws.Range("F25").ClearContents()

# Move the view to where the user was last working (scrolled down, selected J22)
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("J22").Select()
